$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-5, column B currently hold text dates (e.g. "23.03.2024") stored as
# shared strings. The refreshed expense file replaces them with real Excel
# date values (matching the date formatting already used in rows 6-17),
# which also drops the now-unused text-date entries from the shared string
# table and re-packs the remaining label strings.
$ws.Range("B2").Value = [DateTime]"2024-03-23"
$ws.Range("B3").Value = [DateTime]"2024-05-19"
$ws.Range("B4").Value = [DateTime]"2024-05-11"
$ws.Range("B5").Value = [DateTime]"2024-01-07"

# Restore the last selection left behind by the author's edit.
$ws.Range("E6").Select()
